# Applies the edit described by the diff:
#  - Shuffle data rows 4,6,7,9,10,11,12,14,15,16 (rows 2,3,5,8,13 stay put)
#  - Remove the header formatting (bold font / thin border / center-top
#    alignment) so the header cells fall back to the default style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New order of data rows (A:D), keyed by target row number ---------
$rows = @{
    4  = @("mm_2", 381, 4, "/mm/deburr, /mm/drill, /mm/mill, /mm/transport_from_to")
    6  = @("pm_1", 204, 3, "/pm/punch_gill, /pm/punch_recesses, /pm/punch_ribbing")
    7  = @("hw_1", 522, 1, "/hw/human_review")
    9  = @("sm_1", 378, 2, "/sm/sort, /sm/transport")
    10 = @("dm_2", 177, 3, "/dm/cylindrical_drill, /dm/drill, /dm/lower")
    11 = @("wt_1", 447, 1, "/wt/pick_up_and_transport")
    12 = @("ov_2", 330, 1, "/ov/burn")
    14 = @("sm_2", 309, 2, "/sm/sort, /sm/transport")
    15 = @("hbw_2", 1581, 2, "/hbw/store_empty_bucket, /hbw/unload")
    16 = @("vgr_2", 885, 1, "/vgr/pick_up_and_transport")
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
}

# --- Strip the header (row 1) formatting back to the default style -----
$headerRange = $ws.Range("A1:D1")
$headerRange.ClearFormats()
